# Generate Report for Handoff
#
# The localization status moved from "In Translation" to "Ready for
# handoff" and the report was regenerated, so the status text and the
# associated timestamps are refreshed on every sheet. The status column
# is also widened slightly to fit the new, longer label.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# Columns E (zh-cn) and F (de-de) hold the per-locale status; column G
# holds the latest HO Xliff generation timestamp for the row.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 00:37:07"

# --- zh-cn sheet -------------------------------------------------------
# Column C holds Status; column H holds Latest Handoff Datetime.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 00:36:58"

# --- de-de sheet -------------------------------------------------------
# Column C holds Status; column H holds Latest Handoff Datetime.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 00:37:07"

# --- Widen the status columns so the longer text fits ------------------
$wsOverview.Range("E1").ColumnWidth = 16.25
$wsOverview.Range("F1").ColumnWidth = 16.25
$wsZhCn.Range("C1").ColumnWidth = 16.25
$wsDeDe.Range("C1").ColumnWidth = 16.25
